$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear column T (buybacks_cash_returned) for data rows - removed in target
$ws.Range("T2:T3").ClearContents()

# Update capital structure / margin figures for rows 2 and 3
$ws.Range("D2").Value2 = -0.216
$ws.Range("G2").Value2 = -0.1117647058823529
$ws.Range("H2").Value2 = -0.1117647058823529
$ws.Range("I2").Value2 = -0.3642533936651584
$ws.Range("J2").Value2 = -0.3642533936651584
$ws.Range("K2").Value2 = -10.2
$ws.Range("L2").Value2 = -0.4615384615384615
$ws.Range("M2").Value2 = 0
$ws.Range("N2").Value2 = 0
$ws.Range("O2").Value2 = -0
$ws.Range("P2").Value2 = 0
$ws.Range("Q2").Value2 = 0
$ws.Range("R2").Value2 = -0
$ws.Range("U2").Value2 = 36.8
$ws.Range("V2").Value2 = 0.7145631067961165
$ws.Range("W2").Value2 = -0.1522388059701492
$ws.Range("X2").Value2 = 0.1178239745687822
$ws.Range("Y2").Value2 = -0.2700627805389315
$ws.Range("Z2").Value2 = 1.099502487562189
$ws.Range("AA2").Value2 = -0.4004975124378111
$ws.Range("AB2").Value2 = 0.08809070229660761
$ws.Range("AC2").Value2 = -0.4885882147344187
$ws.Range("AD2").Value2 = 28.6
$ws.Range("AE2").Value2 = 0
$ws.Range("AF2").Value2 = 28.6
$ws.Range("AG2").Value2 = -8.199999999999996
$ws.Range("AH2").Value2 = 0.3570536828963796
$ws.Range("AI2").Value2 = 0.3891156462585034
$ws.Range("AJ2").Value2 = -0.1893764434180137
$ws.Range("AK2").Value2 = -0.2234332425068118
$ws.Range("AL2").Value2 = 1.38
$ws.Range("AM2").Value2 = 0.3099999999999998
$ws.Range("AN2").Value2 = -3.376623376623376
$ws.Range("AO2").Value2 = -5.833333333333334
$ws.Range("AP2").Value2 = 0.9681227863046039
$ws.Range("AQ2").Value2 = -25.96774193548389
$ws.Range("D3").Value2 = -0.216
$ws.Range("G3").Value2 = -0.1117647058823529
$ws.Range("H3").Value2 = -0.1117647058823529
$ws.Range("I3").Value2 = -0.3642533936651584
$ws.Range("J3").Value2 = -0.3642533936651584
$ws.Range("K3").Value2 = -10.2
$ws.Range("L3").Value2 = -0.4615384615384615
$ws.Range("M3").Value2 = -0
$ws.Range("N3").Value2 = -0
$ws.Range("O3").Value2 = 0
$ws.Range("P3").Value2 = -0
$ws.Range("Q3").Value2 = -0
$ws.Range("R3").Value2 = 0
$ws.Range("U3").Value2 = 36.8
$ws.Range("V3").Value2 = 0.7145631067961165
$ws.Range("W3").Value2 = -0.1522388059701492
$ws.Range("X3").Value2 = 0.1178239745687822
$ws.Range("Y3").Value2 = -0.2700627805389315
$ws.Range("Z3").Value2 = 1.099502487562189
$ws.Range("AA3").Value2 = -0.4004975124378111
$ws.Range("AB3").Value2 = 0.08809070229660761
$ws.Range("AC3").Value2 = -0.4885882147344187
$ws.Range("AD3").Value2 = 28.6
$ws.Range("AE3").Value2 = 0
$ws.Range("AF3").Value2 = 28.6
$ws.Range("AG3").Value2 = -8.199999999999996
$ws.Range("AH3").Value2 = 0.3570536828963796
$ws.Range("AI3").Value2 = 0.3891156462585034
$ws.Range("AJ3").Value2 = -0.1893764434180137
$ws.Range("AK3").Value2 = -0.2234332425068118
$ws.Range("AL3").Value2 = 1.38
$ws.Range("AM3").Value2 = 0.3099999999999998
$ws.Range("AN3").Value2 = -3.376623376623376
$ws.Range("AO3").Value2 = -5.833333333333334
$ws.Range("AP3").Value2 = 0.9681227863046039
$ws.Range("AQ3").Value2 = -25.96774193548389